$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 67; this shifts the existing rows 67-91
# down to 69-93 and copies formatting (e.g. the date style on column D)
# from the row above, matching native Excel "Insert Rows" behaviour.
$ws.Rows.Item(67).Resize(2).Insert()

# New row 67 data
$ws.Cells.Item(67, 1).Value2 = 9
$ws.Cells.Item(67, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(67, 3).Value2 = "Metropolitana"
$ws.Cells.Item(67, 4).Value2 = 44985
$ws.Cells.Item(67, 5).Value2 = 13
$ws.Cells.Item(67, 6).Value2 = "Fruta"
$ws.Cells.Item(67, 7).Value2 = 100102
$ws.Cells.Item(67, 8).Value2 = "Cítricos"
$ws.Cells.Item(67, 9).Value2 = 100102006
$ws.Cells.Item(67, 10).Value2 = "Pomelo"
$ws.Cells.Item(67, 11).Value2 = "Start Ruby"
$ws.Cells.Item(67, 12).Value2 = "Primera"
$ws.Cells.Item(67, 13).Value2 = 150
$ws.Cells.Item(67, 14).Value2 = 6000
$ws.Cells.Item(67, 15).Value2 = 6000
$ws.Cells.Item(67, 16).Value2 = 6000
$ws.Cells.Item(67, 17).Value2 = "`$/caja 14 kilos"
$ws.Cells.Item(67, 18).Value2 = "Provincia de Quillota"
$ws.Cells.Item(67, 19).Value2 = 429
$ws.Cells.Item(67, 20).Value2 = 14

# New row 68 data
$ws.Cells.Item(68, 1).Value2 = 9
$ws.Cells.Item(68, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(68, 3).Value2 = "Metropolitana"
$ws.Cells.Item(68, 4).Value2 = 44985
$ws.Cells.Item(68, 5).Value2 = 13
$ws.Cells.Item(68, 6).Value2 = "Fruta"
$ws.Cells.Item(68, 7).Value2 = 100102
$ws.Cells.Item(68, 8).Value2 = "Cítricos"
$ws.Cells.Item(68, 9).Value2 = 100102006
$ws.Cells.Item(68, 10).Value2 = "Pomelo"
$ws.Cells.Item(68, 11).Value2 = "Start Ruby"
$ws.Cells.Item(68, 12).Value2 = "Segunda"
$ws.Cells.Item(68, 13).Value2 = 120
$ws.Cells.Item(68, 14).Value2 = 4000
$ws.Cells.Item(68, 15).Value2 = 4000
$ws.Cells.Item(68, 16).Value2 = 4000
$ws.Cells.Item(68, 17).Value2 = "`$/caja 14 kilos"
$ws.Cells.Item(68, 18).Value2 = "Provincia de Quillota"
$ws.Cells.Item(68, 19).Value2 = 286
$ws.Cells.Item(68, 20).Value2 = 14
